$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date value (serial 45175 = 2023-09-06) for every
# data row (rows 2-358). The update bumps that date to serial 45177 (2023-09-08)
# for all of them.
$ws.Range("C2:C358").Value = (Get-Date -Year 2023 -Month 9 -Day 8 -Hour 0 -Minute 0 -Second 0).Date
